# Apply updated crypto price/volume figures (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.108.99'
$ws.Range("E2").Value = '  -1.97%  '

$ws.Range("D3").Value = '2.158.46'
$ws.Range("E3").Value = '  -2.42%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.62'
$ws.Range("E5").Value = '  -1.93%  '

$ws.Range("E6").Value = '  -3.27%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '69.05'
$ws.Range("E7").Value = '  -5.59%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("E9").Value = '  -6.68%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.45'
$ws.Range("E10").Value = '  -9.57%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0906'
$ws.Range("E11").Value = '  -4.70%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.26'
$ws.Range("E12").Value = '  -5.60%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0999'
$ws.Range("E13").Value = '  -3.08%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.63'
$ws.Range("E14").Value = '  -6.04%  '

$ws.Range("D15").Value = '2.477.94'
$ws.Range("E15").Value = '  -2.64%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.18'
$ws.Range("E16").Value = '  -0.46%  '

$ws.Range("D17").Value = '2.139.14'
$ws.Range("E17").Value = '  -3.44%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.780'
$ws.Range("E18").Value = '  -6.55%  '

$ws.Range("D19").Value = '40.891.43'
$ws.Range("E19").Value = '  -2.28%  '

$ws.Range("D20").Value = '0.0₃0988'
$ws.Range("E20").Value = '  -6.95%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '69.63'
$ws.Range("E21").Value = '  -4.59%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.75'
$ws.Range("E22").Value = '  -6.37%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '224.06'
$ws.Range("E23").Value = '  -2.21%  '

$ws.Range("E24").Value = '  -14.25%  '

$ws.Range("E25").Value = '  -0.17%  '

$ws.Range("E26").Value = '  -10.19%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.58'
$ws.Range("E27").Value = '  -9.88%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.42'
$ws.Range("E28").Value = '  -5.19%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.18'
$ws.Range("E29").Value = '  -3.46%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.16'
$ws.Range("E30").Value = '  -1.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.07'
$ws.Range("E31").Value = '  +0.55%  '

$ws.Range("E32").Value = '  -3.39%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.73'
$ws.Range("E33").Value = '  +1.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0750'
$ws.Range("E34").Value = '  -5.43%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.06'
$ws.Range("E35").Value = '  -8.85%  '

$ws.Range("E36").Value = '  -4.12%  '

$ws.Range("E37").Value = '  -8.16%  '

$ws.Range("E38").Value = '  -4.81%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0277'
$ws.Range("E39").Value = '  -7.20%  '

$ws.Range("E40").Value = '  -3.41%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.36'
$ws.Range("E41").Value = '  -17.51%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.27'
$ws.Range("E42").Value = '  -5.97%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '57.83'
$ws.Range("E43").Value = '  -11.87%  '

$ws.Range("E44").Value = '  -5.74%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.20'
$ws.Range("E45").Value = '  -5.55%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0955'
$ws.Range("E46").Value = '  -4.48%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '96.35'
$ws.Range("E47").Value = '  -7.29%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.06'
$ws.Range("E48").Value = '  -4.50%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.10'
$ws.Range("E49").Value = '  -5.36%  '

$ws.Range("E50").Value = '  -9.41%  '

$ws.Range("E51").Value = '  -3.28%  '
